# Adiciona suporte a células do tipo Boolean no importador de planilhas do Excel
#
# Adds a new "BOOLEAN" column (column E) to Sheet1, with a bold header
# (matching the style already used by the NAME/AMOUNT headers) and a
# TRUE boolean value in the data row below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: "BOOLEAN" in row 1, bold (same look as the other
# text headers A1/B1 which use the bold font style).
$ws.Range("E1").Value = "BOOLEAN"
$ws.Range("E1").Font.Bold = $true

# New data cell: a real Boolean value (TRUE) in row 2.
$ws.Range("E2").Value = $true
